# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 299
$wsExhibit.Range("F7").Value  = 12455
$wsExhibit.Range("F10").Value = 22
$wsExhibit.Range("F12").Value = 190
$wsExhibit.Range("F13").Value = 12313
$wsExhibit.Range("F14").Value = 4867
$wsExhibit.Range("F15").Value = 4771
$wsExhibit.Range("F20").Value = 958

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 299
$wsAll.Range("F9").Value  = 12455
$wsAll.Range("F12").Value = 22
$wsAll.Range("F14").Value = 190
$wsAll.Range("F15").Value = 12313
$wsAll.Range("F16").Value = 4867
$wsAll.Range("F17").Value = 4771
$wsAll.Range("F22").Value = 958

$wb.Save()
